# Adding one more test case and Custom Annotation to the framework

$wb = $excel.ActiveWorkbook
$wsRun  = $wb.Worksheets.Item(1)   # RUNMANAGER
$wsData = $wb.Worksheets.Item(2)   # DATA

# ---------------------------------------------------------------------------
# RUNMANAGER sheet: add a new test case row (row 4)
# ---------------------------------------------------------------------------
$wsRun.Range("A4").Value = "validateUpdateMyInfoTest"
$wsRun.Range("B4").Value = "This test is to validate update feature of My Info"
$wsRun.Range("C4").Value = "Yes"

# column B got a little wider to fit the new (longer) description
$wsRun.Columns.Item(2).ColumnWidth = 43.6

# ---------------------------------------------------------------------------
# DATA sheet: add four new columns (firstname, lastname, empid, dropdownvalue)
# used by the Custom Annotation / new test, and a new data row (row 8)
# ---------------------------------------------------------------------------
$wsData.Range("E1").Value = "firstname"
$wsData.Range("F1").Value = "lastname"
$wsData.Range("G1").Value = "empid"
$wsData.Range("H1").Value = "dropdownvalue"

$wsData.Columns.Item(5).ColumnWidth = 8.6
$wsData.Columns.Item(8).ColumnWidth = 14.3

# existing rows 2-7 get blank (quote-prefixed) placeholders in the new columns
for ($r = 2; $r -le 7; $r++) {
    $wsData.Range("E$r").Value = "'"
    $wsData.Range("F$r").Value = "'"
    $wsData.Range("G$r").Value = "'"
    $wsData.Range("H$r").Value = "'"
}

# flip the execute flag (column B) for the existing rows
$wsData.Range("B3").Value = "Yes"
$wsData.Range("B4").Value = "No"
$wsData.Range("B5").Value = "No"
$wsData.Range("B6").Value = "No"
$wsData.Range("B7").Value = "No"

# new row for the added test case
$wsData.Range("A8").Value = "validateUpdateMyInfoTest"
$wsData.Range("B8").Value = "Yes"
$wsData.Range("C8").Value = "Admin"
$wsData.Range("D8").Value = "admin123"
$wsData.Range("E8").Value = "Abhishek"
$wsData.Range("F8").Value = "Singh"
$wsData.Range("G8").Value = "'765507"
$wsData.Range("H8").Value = "Indian"

# ---------------------------------------------------------------------------
# Selections / active sheet: RUNMANAGER becomes active tab, selection D12;
# DATA sheet selection moves to A12 and is no longer the active tab.
# ---------------------------------------------------------------------------
$wsData.Range("A12").Select()
$wsRun.Activate()
$wsRun.Range("D12").Select()
